$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "30.907.55"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.914.04"
$ws.Range("E3").Value = "  +0.94%  "
Set-TextValue ($ws.Range("D4")) "1.001"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue ($ws.Range("D5")) "238.73"
$ws.Range("E5").Value = "  -3.62%  "
$ws.Range("E6").Value = "  +0.07%  "
Set-TextValue ($ws.Range("D7")) "0.4911"
$ws.Range("E7").Value = "  -0.68%  "
Set-TextValue ($ws.Range("D8")) "0.2957"
$ws.Range("E8").Value = "  -0.33%  "
Set-TextValue ($ws.Range("D9")) "0.06758"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "1.901.79"
$ws.Range("E10").Value = "  +0.38%  "
Set-TextValue ($ws.Range("D11")) "17.00"
$ws.Range("E11").Value = "  -1.76%  "
Set-TextValue ($ws.Range("D12")) "0.07299"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue ($ws.Range("D13")) "89.80"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue ($ws.Range("D14")) "5.126"
$ws.Range("E14").Value = "  +0.42%  "
Set-TextValue ($ws.Range("D15")) "0.6694"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "30.876.67"
$ws.Range("E16").Value = "  +0.50%  "
Set-TextValue ($ws.Range("D17")) "0.000007934"
$ws.Range("E17").Value = "  -1.05%  "
Set-TextValue ($ws.Range("D18")) "13.41"
$ws.Range("E18").Value = "  +0.55%  "
Set-TextValue ($ws.Range("D19")) "1.001"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "2.159.78"
$ws.Range("E20").Value = "  +0.99%  "
Set-TextValue ($ws.Range("D21")) "1.000"
$ws.Range("E21").Value = "  -0.07%  "
Set-TextValue ($ws.Range("D22")) "5.140"
$ws.Range("E22").Value = "  +5.69%  "
Set-TextValue ($ws.Range("D23")) "208.08"
$ws.Range("E23").Value = "  +7.13%  "
$ws.Range("E24").Value = "  +1.65%  "
Set-TextValue ($ws.Range("D25")) "9.651"
$ws.Range("E25").Value = "  +1.85%  "
Set-TextValue ($ws.Range("D26")) "157.25"
$ws.Range("E26").Value = "  +1.10%  "
Set-TextValue ($ws.Range("D27")) "18.86"
$ws.Range("E27").Value = "  -2.47%  "
Set-TextValue ($ws.Range("D28")) "1.966"
$ws.Range("E28").Value = "  +1.84%  "
Set-TextValue ($ws.Range("D29")) "1.419"
$ws.Range("E29").Value = "  +1.34%  "
Set-TextValue ($ws.Range("D30")) "4.323"
$ws.Range("E30").Value = "  -1.16%  "
Set-TextValue ($ws.Range("D31")) "0.09151"
$ws.Range("E31").Value = "  +1.36%  "
Set-TextValue ($ws.Range("D32")) "4.045"
$ws.Range("E32").Value = "  -0.23%  "
Set-TextValue ($ws.Range("D33")) "0.05160"
$ws.Range("E33").Value = "  -1.06%  "
Set-TextValue ($ws.Range("D34")) "0.7508"
$ws.Range("E34").Value = "  +0.55%  "
Set-TextValue ($ws.Range("D35")) "1.115"
$ws.Range("E35").Value = "  -1.32%  "
Set-TextValue ($ws.Range("D36")) "2.717"
$ws.Range("E36").Value = "  -0.48%  "
Set-TextValue ($ws.Range("D37")) "0.01840"
$ws.Range("E37").Value = "  -1.35%  "
Set-TextValue ($ws.Range("D38")) "2.732"
$ws.Range("E38").Value = "  +1.82%  "
Set-TextValue ($ws.Range("D39")) "0.9236"
$ws.Range("E39").Value = "  -1.74%  "
Set-TextValue ($ws.Range("D40")) "2.095"
$ws.Range("E40").Value = "  -3.53%  "
Set-TextValue ($ws.Range("D41")) "0.4494"
$ws.Range("E41").Value = "  +0.86%  "
Set-TextValue ($ws.Range("D42")) "106.87"
$ws.Range("E42").Value = "  +0.36%  "
Set-TextValue ($ws.Range("D43")) "5.863"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("E44").Value = "  +0.65%  "
Set-TextValue ($ws.Range("D45")) "7.713"
$ws.Range("E45").Value = "  -0.07%  "
Set-TextValue ($ws.Range("D46")) "0.1383"
$ws.Range("E46").Value = "  +2.54%  "
Set-TextValue ($ws.Range("D47")) "66.34"
$ws.Range("E47").Value = "  +13.94%  "
Set-TextValue ($ws.Range("D48")) "35.15"
$ws.Range("E48").Value = "  +4.28%  "
Set-TextValue ($ws.Range("D49")) "8.975"
$ws.Range("E49").Value = "  +2.97%  "
Set-TextValue ($ws.Range("D50")) "0.4087"
$ws.Range("E50").Value = "  +3.20%  "
Set-TextValue ($ws.Range("D51")) "0.05920"
$ws.Range("E51").Value = "  +1.03%  "
